$wb = $excel.ActiveWorkbook

# --- References to the existing sheets ---------------------------------
$summary = $wb.Worksheets.Item(1)         # "总计"
$q3      = $wb.Worksheets.Item("2022-Q3") # existing quarter sheet (stays "2022-Q3")

# --- 1) Insert the new "2022-Q4" sheet ----------------------------------
# Duplicate the "总计" sheet (so the new sheet inherits the same sheet
# properties / page margins) right after itself, then rename it and wipe
# its copied data. Copy() activates the freshly created sheet, which is
# how we get a handle on it.
$summary.Copy($null, $summary)
$q4 = $wb.ActiveSheet
$q4.Name = "2022-Q4"
$q4.Cells.ClearContents()

# --- 2) Populate the new "2022-Q4" sheet with the fund snapshot --------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# The source figures are kept as plain text (e.g. the fund code has a
# leading zero), so force the Text number format before entering them -
# otherwise Excel would helpfully "fix" them into numbers.
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "006105"
$q4.Range("C2").Value = "泰达宏利印度机会股票（QDII）"
$q4.Range("D2").Value = "0.67"
$q4.Range("E2").Value = "91.81"
$q4.Range("F2").Value = "2.57"
$q4.Range("G2").Value = "0.0172"
$q4.Range("H2").Value = 9

# Extend the header formatting (bold font + border + centered alignment -
# style index 2 on the "总计" sheet) across the full header row, since the
# copy only brought that style along for columns B:D. (A2 already kept its
# style from the sheet copy, so it needs no extra work.)
$q4.Range("B1:D1").Copy()
$q4.Range("E1:H1").PasteSpecial(-4122) # xlPasteFormats

# --- 3) Update the "总计" (summary) sheet -------------------------------
# Row 2 now describes the newly added quarter.
$summary.Range("B2").Value = "2022-Q4"

# Row 3 is a new row re-stating the previous "2022-Q3" summary figures.
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.02

# Give the new A3 cell the same formatting as A2 (bold, bordered, centered).
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122) # xlPasteFormats
